$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit re-sorts the occurrence rows (2-19) into a new order: each row
# keeps its position but takes on the Id/species/coordinate data that, before
# the edit, belonged to a different row in the same table (row 14 is already
# in its correct place and is left untouched). The writes below set each
# touched cell to the value it has after that re-sort.

# --- Id (A), taxon block (B,D,E,F,G,H) and coordinates (Q,R) ---
$ws.Range("A2").Value2 = 111396323
$ws.Range("B2").Value2 = 96348
$ws.Range("D2").Value2 = "VU"
$ws.Range("E2").Value2 = 220787
$ws.Range("F2").Value2 = "Knärot"
$ws.Range("G2").Value2 = "Goodyera repens"
$ws.Range("H2").Value2 = "(L.) R. Br."
$ws.Range("Q2").Value2 = 625301.6605433678
$ws.Range("R2").Value2 = 7209610.70454926
$ws.Range("A3").Value2 = 111396309
$ws.Range("Q3").Value2 = 625341.71034419
$ws.Range("R3").Value2 = 7209536.108963673
$ws.Range("A4").Value2 = 111396312
$ws.Range("Q4").Value2 = 625242.7087276473
$ws.Range("R4").Value2 = 7209468.80281719
$ws.Range("A5").Value2 = 111396311
$ws.Range("Q5").Value2 = 625271.0561409625
$ws.Range("R5").Value2 = 7209511.101565193
$ws.Range("A6").Value2 = 111396324
$ws.Range("Q6").Value2 = 625335.6676841485
$ws.Range("R6").Value2 = 7209609.168182318
$ws.Range("A7").Value2 = 111396325
$ws.Range("Q7").Value2 = 625389.9085714296
$ws.Range("R7").Value2 = 7209580.514361567
$ws.Range("A8").Value2 = 111396313
$ws.Range("Q8").Value2 = 625231.5510770321
$ws.Range("R8").Value2 = 7209481.895207534
$ws.Range("A9").Value2 = 111396318
$ws.Range("Q9").Value2 = 625177.6865340136
$ws.Range("R9").Value2 = 7209552.099144561
$ws.Range("A10").Value2 = 111396319
$ws.Range("Q10").Value2 = 625228.8129008666
$ws.Range("R10").Value2 = 7209607.642547456
$ws.Range("A11").Value2 = 111396316
$ws.Range("Q11").Value2 = 625153.7279882778
$ws.Range("R11").Value2 = 7209526.513740451
$ws.Range("A12").Value2 = 111396310
$ws.Range("Q12").Value2 = 625289.0018867656
$ws.Range("R12").Value2 = 7209518.212698339
$ws.Range("A13").Value2 = 111396317
$ws.Range("Q13").Value2 = 625153.5624699651
$ws.Range("R13").Value2 = 7209550.662191558
$ws.Range("A15").Value2 = 111396315
$ws.Range("Q15").Value2 = 625167.9685939638
$ws.Range("R15").Value2 = 7209530.9258211
$ws.Range("A16").Value2 = 111396321
$ws.Range("Q16").Value2 = 625240.2002264742
$ws.Range("R16").Value2 = 7209649.650274927
$ws.Range("A17").Value2 = 111396326
$ws.Range("Q17").Value2 = 625397.1584455093
$ws.Range("R17").Value2 = 7209589.718691397
$ws.Range("A18").Value2 = 111396308
$ws.Range("B18").Value2 = 56398
$ws.Range("D18").Value2 = "NT"
$ws.Range("E18").Value2 = 100109
$ws.Range("F18").Value2 = "Tretåig hackspett"
$ws.Range("G18").Value2 = "Picoides tridactylus"
$ws.Range("H18").Value2 = "(Linnaeus, 1758)"
$ws.Range("M18").Value2 = "äldre spår"
$ws.Range("Q18").Value2 = 625151.1577179903
$ws.Range("R18").Value2 = 7209567.512248591
$ws.Range("A19").Value2 = 111396314
$ws.Range("Q19").Value2 = 625202.8383709632
$ws.Range("R19").Value2 = 7209539.171001118

# --- row 2 loses its "aldre spar" activity-note block (K,L,M,N), which
#     only row 2 used to carry ---
$ws.Range("K2").Value2 = $null
$ws.Range("L2").Value2 = $null
$ws.Range("M2").Value2 = $null
$ws.Range("N2").Value2 = $null

# --- row 18 gains that activity-note block; K18/L18/N18 stay blank (only
#     M18 carries text) but must exist as real (empty) cells, so force cell
#     creation via a text format before writing the blank value ---
$ws.Range("K18").NumberFormat = "@"
$ws.Range("K18").Value2 = ""
$ws.Range("K18").Style = "Normal"
$ws.Range("L18").NumberFormat = "@"
$ws.Range("L18").Value2 = ""
$ws.Range("L18").Style = "Normal"
$ws.Range("N18").NumberFormat = "@"
$ws.Range("N18").Value2 = ""
$ws.Range("N18").Style = "Normal"
$ws.Range("M18").Value2 = "äldre spår"
